$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E2:E12 values from 50 to 70
$ws.Range("E2:E12").Value = 70

# Update the selection to reflect the new active cell (E17, single cell)
$ws.Range("E17").Select()
